$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Step 1 (OPEN_URL) now reported as a FAIL with a null-reference error ---
$ws.Range("L2").Value = "FAIL"
$ws.Range("M2").Value = "Cannot read properties of null (reading 'isClosed')"
$ws.Range("N2").Value = "Cannot read properties of null (reading 'isClosed')"

# O2/P2 become empty-string (Text) results rather than screenshot/page-source paths.
# A plain Value = "" clears the cell entirely, so force an empty text entry via the
# leading-apostrophe trick, then drop the resulting "quote prefix" style so the cell
# matches a normal, unstyled, empty text cell.
$ws.Range("O2").Value = "'"
$ws.Range("O2").Style = "Normal"
$ws.Range("P2").Value = "'"
$ws.Range("P2").Style = "Normal"

# --- Rows 3-6: these steps are no longer executed/reported at all -> drop Status/Remarks ---
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()

$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()

$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()

$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()

# --- Rows 7-9: these steps' result columns (Status/Remarks/Actual Output/Screenshot/Page Source) are removed ---
$ws.Range("L7:P7").ClearContents()
$ws.Range("L8:P8").ClearContents()
$ws.Range("L9:P9").ClearContents()
